$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the style of the existing
# header cells (e.g. G1) by copying G1's formatting over to H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-5
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
